# Workbook "Hortaliza, Feria Lagunitas de Puerto Montt - Cebolla"
# Weekly fruit/vegetable price update: two new observations were recorded
# for this market/product (dated 44939) and inserted at the top of the
# existing data block (rows 744-745), pushing the previously-existing
# rows 744:773 down to 746:775.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at the top of the data block, shifting the
# existing rows (old 744:773) down to 746:775.
$ws.Rows("744:745").Insert()

# --- New row 744 ---
$ws.Cells.Item(744, 1).Value = 4
$ws.Cells.Item(744, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(744, 3).Value = "Los Lagos"
$ws.Cells.Item(744, 4).Value = 44939
$ws.Cells.Item(744, 5).Value = 10
$ws.Cells.Item(744, 6).Value = 100112004
$ws.Cells.Item(744, 7).Value = "Cebolla"
$ws.Cells.Item(744, 8).Value = "Morada(o)"
$ws.Cells.Item(744, 9).Value = "1a (cosecha)"
$ws.Cells.Item(744, 10).Value = 250
$ws.Cells.Item(744, 11).Value = 16000
$ws.Cells.Item(744, 12).Value = 16000
$ws.Cells.Item(744, 13).Value = 16000
$ws.Cells.Item(744, 14).Value = "`$/malla 18 kilos"
$ws.Cells.Item(744, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(744, 16).Value = 889
$ws.Cells.Item(744, 17).Value = 18
$ws.Cells.Item(744, 18).Value = "Hortaliza"

# --- New row 745 ---
$ws.Cells.Item(745, 1).Value = 4
$ws.Cells.Item(745, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(745, 3).Value = "Los Lagos"
$ws.Cells.Item(745, 4).Value = 44939
$ws.Cells.Item(745, 5).Value = 10
$ws.Cells.Item(745, 6).Value = 100112004
$ws.Cells.Item(745, 7).Value = "Cebolla"
$ws.Cells.Item(745, 8).Value = "Sin especificar"
$ws.Cells.Item(745, 9).Value = "1a (cosecha)"
$ws.Cells.Item(745, 10).Value = 1000
$ws.Cells.Item(745, 11).Value = 15000
$ws.Cells.Item(745, 12).Value = 15000
$ws.Cells.Item(745, 13).Value = 15000
$ws.Cells.Item(745, 14).Value = "`$/malla 18 kilos"
$ws.Cells.Item(745, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(745, 16).Value = 833
$ws.Cells.Item(745, 17).Value = 18
$ws.Cells.Item(745, 18).Value = "Hortaliza"
